$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.118.08'
$ws.Range("E2").Value = '  -2.41%  '

$ws.Range("D3").Value = '1.646.72'
$ws.Range("E3").Value = '  -1.96%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.48%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3919'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3867'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.56%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.003'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.370'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '49.41'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08574'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.61'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.100'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.43%  '

$ws.Range("E15").Value = '  -2.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.521'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.77%  '

$ws.Range("D17").Value = '1.648.31'
$ws.Range("E17").Value = '  -5.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06909'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.908'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.82%  '

$ws.Range("D24").Value = '24.119.07'
$ws.Range("E24").Value = '  -2.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.436'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.884'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '158.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.76%  '

$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.244'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.67%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '140.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.305'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.492'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.10%  '

$ws.Range("D33").Value = '1.820.62'
$ws.Range("E33").Value = '  -1.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08160'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.819'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02914'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9671'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2690'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09174'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.36'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.426'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7529'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6908'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.463'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.105'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.56%  '

$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08379'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.268'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.90%  '
